# Auto-generated edit script: update crypto price/volume table to new snapshot values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.020.40'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.619.75'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.37'
$ws.Range("E5").Value = '  -1.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.512'
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.250'
$ws.Range("E8").Value = '  -1.28%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0627'
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.93'
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("D12").Value = '1.847.31'
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").Value = '1.620.27'
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.11'
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.537'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").Value = '27.015.86'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.40'
$ws.Range("E17").Value = '  -3.25%  '
$ws.Range("D18").Value = '0.0₃0736'
$ws.Range("E18").Value = '  -0.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.36'
$ws.Range("E19").Value = '  -1.55%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.82'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.32'
$ws.Range("E22").Value = '  -2.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.36'
$ws.Range("E23").Value = '  -6.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.97'
$ws.Range("E24").Value = '  -1.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.39'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.47'
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("B27").Value = 'BinanceUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E28").Value = '  -3.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.50'
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0509'
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("E32").Value = '  -2.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.701'
$ws.Range("E33").Value = '  +28.74%  '
$ws.Range("E34").Value = '  -1.14%  '
$ws.Range("D35").Value = '1.334.90'
$ws.Range("E35").Value = '  +2.61%  '
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0175'
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.838'
$ws.Range("E39").Value = '  -1.71%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.797'
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("E42").Value = '  -0.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.35'
$ws.Range("E43").Value = '  +0.59%  '
$ws.Range("E44").Value = '  +3.54%  '
$ws.Range("D45").Value = '1.758.62'
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.90'
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.63'
$ws.Range("E47").Value = '  +2.21%  '
$ws.Range("E48").Value = '  +25.16%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0513'
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0993'
$ws.Range("E50").Value = '  +3.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.56'
$ws.Range("E51").Value = '  -0.62%  '
